# Lazy load access rights (#368) - add "Send Payment Notification" and
# "Send Call Notice" columns to the CapitalCall sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCall")

# New header row (row 1), columns L and M
$ws.Cells.Item(1, 12).Value = "Send Payment Notification"
$ws.Cells.Item(1, 13).Value = "Send Call Notice"

# Default both new flags to "Yes" for every data row (2-4)
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 12).Value = "Yes"
    $ws.Cells.Item($r, 13).Value = "Yes"
}

# Match the saved selection state from the authored workbook
$ws.Range("K1").Select()
